$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3.0
$ws.Range("G2").Value = 216.5832213333333
$ws.Range("H2").Value = 649.749664
$ws.Range("I2").Value = 0.4331411212367192
$ws.Range("J2").Value = 0.4331411212367192
$ws.Range("K2").Value = 3.0
$ws.Range("M2").Value = 135.955556
$ws.Range("N2").Value = 407.866668
$ws.Range("O2").Value = 0.6947679994035034
$ws.Range("P2").Value = 0.6947679994035034
$ws.Range("Q2").Value = 29445.6922766444
$ws.Range("R2").Value = 265011.2304897996
$ws.Range("S2").Value = 0.3009325902610258
$ws.Range("T2").Value = 0.3009325902610257

# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3.0
$ws.Range("G3").Value = 216.5832213333333
$ws.Range("H3").Value = 649.749664
$ws.Range("I3").Value = 0.4331411212367192
$ws.Range("J3").Value = 0.4331411212367192
$ws.Range("K3").Value = 2.0
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.449122
$ws.Range("N3").Value = 1.347366
$ws.Range("O3").Value = 0.002295129398228494
$ws.Range("P3").Value = 0.002295129398228494
$ws.Range("Q3").Value = 97.27228953166934
$ws.Range("R3").Value = 875.4506057850241
$ws.Range("S3").Value = 0.0009941149209320465
$ws.Range("T3").Value = 0.0009941149209320463

# Row 4
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3.0
$ws.Range("G4").Value = 216.5832213333333
$ws.Range("H4").Value = 649.749664
$ws.Range("I4").Value = 0.4331411212367192
$ws.Range("J4").Value = 0.4331411212367192
$ws.Range("K4").Value = 3.0
$ws.Range("M4").Value = 56.38366533333333
$ws.Range("N4").Value = 169.150996
$ws.Range("O4").Value = 0.2881350899898248
$ws.Range("P4").Value = 0.2881350899898248
$ws.Range("Q4").Value = 12211.75586847393
$ws.Range("R4").Value = 109905.8028162653
$ws.Range("S4").Value = 0.1248031559458357
$ws.Range("T4").Value = 0.1248031559458357

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3.0
$ws.Range("G5").Value = 216.5832213333333
$ws.Range("H5").Value = 649.749664
$ws.Range("I5").Value = 0.4331411212367192
$ws.Range("J5").Value = 0.4331411212367192
$ws.Range("K5").Value = 3.0
$ws.Range("M5").Value = 2.896484
$ws.Range("N5").Value = 8.689452
$ws.Range("O5").Value = 0.01480178120844327
$ws.Range("P5").Value = 0.01480178120844327
$ws.Range("Q5").Value = 627.3298352604586
$ws.Range("R5").Value = 5645.968517344128
$ws.Range("S5").Value = 0.00641126010892572
$ws.Range("T5").Value = 0.006411260108925719

# Row 6
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3.0
$ws.Range("G6").Value = 161.954974
$ws.Range("H6").Value = 485.864922
$ws.Range("I6").Value = 0.3238910133313607
$ws.Range("J6").Value = 0.3238910133313606
$ws.Range("K6").Value = 3.0
$ws.Range("M6").Value = 135.955556
$ws.Range("N6").Value = 407.866668
$ws.Range("O6").Value = 0.6947679994035034
$ws.Range("P6").Value = 0.6947679994035034
$ws.Range("Q6").Value = 22018.67853713554
$ws.Range("R6").Value = 198168.1068342199
$ws.Range("S6").Value = 0.2250291113570029
$ws.Range("T6").Value = 0.2250291113570029

# Row 7
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3.0
$ws.Range("G7").Value = 161.954974
$ws.Range("H7").Value = 485.864922
$ws.Range("I7").Value = 0.3238910133313607
$ws.Range("J7").Value = 0.3238910133313606
$ws.Range("K7").Value = 2.0
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.449122
$ws.Range("N7").Value = 1.347366
$ws.Range("O7").Value = 0.002295129398228494
$ws.Range("P7").Value = 0.002295129398228494
$ws.Range("Q7").Value = 72.737541832828
$ws.Range("R7").Value = 654.637876495452
$ws.Range("S7").Value = 0.0007433717865188229
$ws.Range("T7").Value = 0.0007433717865188227

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3.0
$ws.Range("G8").Value = 161.954974
$ws.Range("H8").Value = 485.864922
$ws.Range("I8").Value = 0.3238910133313607
$ws.Range("J8").Value = 0.3238910133313606
$ws.Range("K8").Value = 3.0
$ws.Range("M8").Value = 56.38366533333333
$ws.Range("N8").Value = 169.150996
$ws.Range("O8").Value = 0.2881350899898248
$ws.Range("P8").Value = 0.2881350899898248
$ws.Range("Q8").Value = 9131.6150530847
$ws.Range("R8").Value = 82184.5354777623
$ws.Range("S8").Value = 0.09332436627312717
$ws.Range("T8").Value = 0.09332436627312712

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3.0
$ws.Range("G9").Value = 161.954974
$ws.Range("H9").Value = 485.864922
$ws.Range("I9").Value = 0.3238910133313607
$ws.Range("J9").Value = 0.3238910133313606
$ws.Range("K9").Value = 3.0
$ws.Range("M9").Value = 2.896484
$ws.Range("N9").Value = 8.689452
$ws.Range("O9").Value = 0.01480178120844327
$ws.Range("P9").Value = 0.01480178120844327
$ws.Range("Q9").Value = 469.0999909114159
$ws.Range("R9").Value = 4221.899918202744
$ws.Range("S9").Value = 0.004794163914711784
$ws.Range("T9").Value = 0.004794163914711783

# Row 10
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3.0
$ws.Range("G10").Value = 0.4608033333333333
$ws.Range("H10").Value = 1.38241
$ws.Range("I10").Value = 0.0009215527926904059
$ws.Range("J10").Value = 0.0009215527926904059
$ws.Range("K10").Value = 3.0
$ws.Range("M10").Value = 135.955556
$ws.Range("N10").Value = 407.866668
$ws.Range("O10").Value = 0.6947679994035034
$ws.Range("P10").Value = 0.6947679994035034
$ws.Range("Q10").Value = 62.64877338998666
$ws.Range("R10").Value = 563.8389605098799
$ws.Range("S10").Value = 0.0006402653901222249
$ws.Range("T10").Value = 0.0006402653901222249

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3.0
$ws.Range("G11").Value = 0.4608033333333333
$ws.Range("H11").Value = 1.38241
$ws.Range("I11").Value = 0.0009215527926904059
$ws.Range("J11").Value = 0.0009215527926904059
$ws.Range("K11").Value = 2.0
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.449122
$ws.Range("N11").Value = 1.347366
$ws.Range("O11").Value = 0.002295129398228494
$ws.Range("P11").Value = 0.002295129398228494
$ws.Range("Q11").Value = 0.2069569146733333
$ws.Range("R11").Value = 1.86261223206
$ws.Range("S11").Value = 0.000002115082906523319
$ws.Range("T11").Value = 0.000002115082906523319

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("E12").Value = 3.0
$ws.Range("G12").Value = 0.4608033333333333
$ws.Range("H12").Value = 1.38241
$ws.Range("I12").Value = 0.0009215527926904059
$ws.Range("J12").Value = 0.0009215527926904059
$ws.Range("K12").Value = 3.0
$ws.Range("M12").Value = 56.38366533333333
$ws.Range("N12").Value = 169.150996
$ws.Range("O12").Value = 0.2881350899898248
$ws.Range("P12").Value = 0.2881350899898248
$ws.Range("Q12").Value = 25.98178093115111
$ws.Range("R12").Value = 233.83602838036
$ws.Range("S12").Value = 0.0002655316968522245
$ws.Range("T12").Value = 0.0002655316968522245

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("E13").Value = 3.0
$ws.Range("G13").Value = 0.4608033333333333
$ws.Range("H13").Value = 1.38241
$ws.Range("I13").Value = 0.0009215527926904059
$ws.Range("J13").Value = 0.0009215527926904059
$ws.Range("K13").Value = 3.0
$ws.Range("M13").Value = 2.896484
$ws.Range("N13").Value = 8.689452
$ws.Range("O13").Value = 0.01480178120844327
$ws.Range("P13").Value = 0.01480178120844327
$ws.Range("Q13").Value = 1.334709482146666
$ws.Range("R13").Value = 12.01238533932
$ws.Range("S13").Value = 0.00001364062280943327
$ws.Range("T13").Value = 0.00001364062280943327

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Col4a1"
$ws.Range("C14").Value = "Cd93"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 121.0302313333333
$ws.Range("H14").Value = 363.090694
$ws.Range("I14").Value = 0.2420463126392298
$ws.Range("J14").Value = 0.2420463126392298
$ws.Range("K14").Value = 3.0
$ws.Range("L14").Value = 1.0
$ws.Range("M14").Value = 135.955556
$ws.Range("N14").Value = 407.866668
$ws.Range("O14").Value = 0.6947679994035034
$ws.Range("P14").Value = 0.6947679994035034
$ws.Range("Q14").Value = 16454.73239373196
$ws.Range("R14").Value = 148092.5915435876
$ws.Range("S14").Value = 0.1681660323953526
$ws.Range("T14").Value = 0.1681660323953526

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Col4a1"
$ws.Range("C15").Value = "Cd93"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 121.0302313333333
$ws.Range("H15").Value = 363.090694
$ws.Range("I15").Value = 0.2420463126392298
$ws.Range("J15").Value = 0.2420463126392298
$ws.Range("K15").Value = 2.0
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.449122
$ws.Range("N15").Value = 1.347366
$ws.Range("O15").Value = 0.002295129398228494
$ws.Range("P15").Value = 0.002295129398228494
$ws.Range("Q15").Value = 54.35733955688934
$ws.Range("R15").Value = 489.216056012004
$ws.Range("S15").Value = 0.0005555276078711014
$ws.Range("T15").Value = 0.0005555276078711012

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Col4a1"
$ws.Range("C16").Value = "Cd93"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 121.0302313333333
$ws.Range("H16").Value = 363.090694
$ws.Range("I16").Value = 0.2420463126392298
$ws.Range("J16").Value = 0.2420463126392298
$ws.Range("K16").Value = 3.0
$ws.Range("L16").Value = 1.0
$ws.Range("M16").Value = 56.38366533333333
$ws.Range("N16").Value = 169.150996
$ws.Range("O16").Value = 0.2881350899898248
$ws.Range("P16").Value = 0.2881350899898248
$ws.Range("Q16").Value = 6824.12805871458
$ws.Range("R16").Value = 61417.15252843122
$ws.Range("S16").Value = 0.06974203607400975
$ws.Range("T16").Value = 0.06974203607400974

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Col4a1"
$ws.Range("C17").Value = "Cd93"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3.0
$ws.Range("F17").Value = 1.0
$ws.Range("G17").Value = 121.0302313333333
$ws.Range("H17").Value = 363.090694
$ws.Range("I17").Value = 0.2420463126392298
$ws.Range("J17").Value = 0.2420463126392298
$ws.Range("K17").Value = 3.0
$ws.Range("L17").Value = 1.0
$ws.Range("M17").Value = 2.896484
$ws.Range("N17").Value = 8.689452
$ws.Range("O17").Value = 0.01480178120844327
$ws.Range("P17").Value = 0.01480178120844327
$ws.Range("Q17").Value = 350.5621285732986
$ws.Range("R17").Value = 3155.059157159687
$ws.Range("S17").Value = 0.003582716561996337
$ws.Range("T17").Value = 0.003582716561996337
